{"js": "const body = context.document.body;\n\n// --- Change 1 (typo fix, lecture-note text): \"\u043f\u0435\u0440\u0435\u0441\u043a\u0430\u0437\u044b\u0432\u0430\u0442\u044c\" -> \"\u043f\u0440\u0435\u0434\u0441\u043a\u0430\u0437\u044b\u0432\u0430\u0442\u044c\" ---\nconst typoResults = body.search(\"\u043f\u0435\u0440\u0435\u0441\u043a\u0430\u0437\u044b\u0432\u0430\u0442\u044c \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u0439 \u0431\u0438\u0442\", { matchCase: true });\ntypoResults.load(\"items\");\nawait context.sync();\nif (typoResults.items.length > 0) {\n  typoResults.items[0].insertText(\"\u043f\u0440\u0435\u0434\u0441\u043a\u0430\u0437\u044b\u0432\u0430\u0442\u044c \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u0439 \u0431\u0438\u0442\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Change 2: the cursor's last-edit marker (\"_GoBack\") now sits inside\n// \"..., \u043f\u043e\u0437\u0432\u043e\u043b\u044f\u044e\u0449\u0438\u0439 \u0432\u044b|\u0438\u0433\u0440\u0430\u0442\u044c \u0438\u0433\u0440\u0443...\" (this is where the author's cursor\n// was left after their most recent edit). Move/insert the bookmark there. ---\nconst gameResults = body.search(\", \u043f\u043e\u0437\u0432\u043e\u043b\u044f\u044e\u0449\u0438\u0439 \u0432\u044b\", { matchCase: true });\ngameResults.load(\"items\");\nawait context.sync();\nif (gameResults.items.length > 0) {\n  const splitPoint = gameResults.items[0].getRange(Word.RangeLocation.end);\n  splitPoint.insertBookmark(\"_GoBack\");\n}\nawait context.sync();\n\n// --- Change 3: \"... \u043f\u0440\u043e\u0434\u0435\u043c\u043e\u043d\u0441\u0442\u0440\u0438\u0440\u0443\u0439\u0442\u0435\" / (old \"_GoBack\") / \" \u0430\u0442\u0430\u043a\u0443.\" is\n// simplified back into one contiguous run of text (the stray bookmark that\n// used to interrupt it is removed by this replace). ---\nconst attackResults = body.search(\" \u043f\u0440\u043e\u0434\u0435\u043c\u043e\u043d\u0441\u0442\u0440\u0438\u0440\u0443\u0439\u0442\u0435 \u0430\u0442\u0430\u043a\u0443.\", { matchCase: true });\nattackResults.load(\"items\");\nawait context.sync();\nif (attackResults.items.length > 0) {\n  attackResults.items[0].insertText(\" \u043f\u0440\u043e\u0434\u0435\u043c\u043e\u043d\u0441\u0442\u0440\u0438\u0440\u0443\u0439\u0442\u0435 \u0430\u0442\u0430\u043a\u0443.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1 (typo fix): \"\u043f\u0435\u0440\u0435\u0441\u043a\u0430\u0437\u044b\u0432\u0430\u0442\u044c\" -> \"\u043f\u0440\u0435\u0434\u0441\u043a\u0430\u0437\u044b\u0432\u0430\u0442\u044c\" ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"\u043f\u0435\u0440\u0435\u0441\u043a\u0430\u0437\u044b\u0432\u0430\u0442\u044c \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u0439 \u0431\u0438\u0442\"\n$find1.Replacement.Text = \"\u043f\u0440\u0435\u0434\u0441\u043a\u0430\u0437\u044b\u0432\u0430\u0442\u044c \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u0439 \u0431\u0438\u0442\"\n$find1.Execute([ref]$find1.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$find1.Replacement.Text, [ref]2)\n\n# --- Change 3 (do this before change 2 so the old bookmark -- not the\n# relocated one -- is what gets removed here): \"... \u043f\u0440\u043e\u0434\u0435\u043c\u043e\u043d\u0441\u0442\u0440\u0438\u0440\u0443\u0439\u0442\u0435\" /\n# (old bookmark) / \" \u0430\u0442\u0430\u043a\u0443.\" becomes one contiguous run of text again. ---\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"\u0430\u0442\u0430\u043a\u0443.\"\n$find3.Replacement.Text = \"\u0430\u0442\u0430\u043a\u0443.\"\n$find3.Execute([ref]$find3.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$find3.Replacement.Text, [ref]2)\n\n# --- Change 2: move the \"_GoBack\" (last-edit-location) bookmark into\n# \"..., \u043f\u043e\u0437\u0432\u043e\u043b\u044f\u044e\u0449\u0438\u0439 \u0432\u044b|\u0438\u0433\u0440\u0430\u0442\u044c \u0438\u0433\u0440\u0443...\" ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \", \u043f\u043e\u0437\u0432\u043e\u043b\u044f\u044e\u0449\u0438\u0439 \u0432\u044b\"\n$find2.Execute()\nif ($find2.Found) {\n    $splitPoint = $find2.Parent.Duplicate\n    $splitPoint.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $splitPoint)\n}\n"}
